$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 4901
$ws.Range("L3").Value = 5280
$ws.Range("D4").Value = 1997
$ws.Range("F4").Value = 1932
$ws.Range("H4").Value = 1763
$ws.Range("L4").Value = 1290
$ws.Range("L5").Value = 311
$ws.Range("L6").Value = 4452
$ws.Range("D7").Value = 28188
$ws.Range("F7").Value = 24125
$ws.Range("H7").Value = 26079
$ws.Range("L7").Value = 16234

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L7").Value = 530
$ws.Range("L8").Value = 1081
$ws.Range("L11").Value = 263
$ws.Range("L18").Value = 115
$ws.Range("L19").Value = 445
$ws.Range("L20").Value = 407
$ws.Range("L23").Value = 177
$ws.Range("L25").Value = 96
$ws.Range("L29").Value = 886
$ws.Range("L33").Value = 745
$ws.Range("L35").Value = 24
$ws.Range("L36").Value = 212
$ws.Range("L37").Value = 613
$ws.Range("L39").Value = 10
$ws.Range("L42").Value = 528
$ws.Range("L46").Value = 36
$ws.Range("L49").Value = 82
$ws.Range("L50").Value = 82
$ws.Range("L51").Value = 206
$ws.Range("L52").Value = 328
$ws.Range("L54").Value = 343
$ws.Range("L60").Value = 104
$ws.Range("D63").Value = 377
$ws.Range("F63").Value = 216
$ws.Range("H63").Value = 314
$ws.Range("L63").Value = 48
$ws.Range("L65").Value = 319
$ws.Range("L66").Value = 41
$ws.Range("L71").Value = 46
$ws.Range("L73").Value = 127
$ws.Range("L82").Value = 24
$ws.Range("L84").Value = 159
$ws.Range("L85").Value = 830
$ws.Range("L88").Value = 175
$ws.Range("L94").Value = 199
$ws.Range("L95").Value = 224
$ws.Range("L97").Value = 138
$ws.Range("L99").Value = 279
$ws.Range("D101").Value = 28188
$ws.Range("F101").Value = 24125
$ws.Range("H101").Value = 26079
$ws.Range("L101").Value = 16234

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 182
$ws.Range("L7").Value = 530

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 99
$ws.Range("L3").Value = 80
$ws.Range("L7").Value = 263

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 338
$ws.Range("L7").Value = 830

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 106
$ws.Range("L3").Value = 104
$ws.Range("L7").Value = 328

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L3").Value = 368
$ws.Range("L6").Value = 279
$ws.Range("L7").Value = 1081

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 259
$ws.Range("L6").Value = 223
$ws.Range("L7").Value = 745

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L2").Value = 88
$ws.Range("L5").Value = 7
$ws.Range("L7").Value = 224

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 184
$ws.Range("L3").Value = 207
$ws.Range("L7").Value = 613

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 115
$ws.Range("L7").Value = 319

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L3").Value = 120
$ws.Range("L7").Value = 279

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L2").Value = 52
$ws.Range("L7").Value = 159

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 82

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L6").Value = 167
$ws.Range("L7").Value = 343

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 336
$ws.Range("L7").Value = 886

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 157
$ws.Range("L3").Value = 139
$ws.Range("L4").Value = 17
$ws.Range("L6").Value = 125
$ws.Range("L7").Value = 445

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L6").Value = 147
$ws.Range("L7").Value = 528

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("L3").Value = 11
$ws.Range("L7").Value = 36

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L6").Value = 48
$ws.Range("L7").Value = 177

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 133
$ws.Range("L7").Value = 407

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L2").Value = 42
$ws.Range("L7").Value = 115

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L3").Value = 64
$ws.Range("L7").Value = 212

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L3").Value = 46
$ws.Range("L5").Value = 5
$ws.Range("L6").Value = 79
$ws.Range("L7").Value = 199

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L5").Value = 3
$ws.Range("L7").Value = 96

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L6").Value = 21
$ws.Range("L7").Value = 82

$ws = $wb.Worksheets.Item("Greektown")
$ws.Range("L2").Value = 4
$ws.Range("L6").Value = 10

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("L3").Value = 11
$ws.Range("L7").Value = 41

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 24

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 127

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L2").Value = 33
$ws.Range("L7").Value = 138

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L4").Value = 8
$ws.Range("L7").Value = 175

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L2").Value = 62
$ws.Range("L7").Value = 206

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 104

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("L2").Value = 20
$ws.Range("L7").Value = 46

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("L3").Value = 8
$ws.Range("L7").Value = 24
